# Update the "_old" baseline model columns (Ada_old, Avey_old, Buoy_old,
# K health_old, WebMD_old, doctor_MA_old, doctor_NJ_old, doctor_TH_old)
# on Sheet1 with recomputed metrics after adding the Harvard case
# classification to the gold-standard set (precision, recall, f1-score,
# f2-score, NDCG and length-vs-gs rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("I2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.4
$ws.Range("Q2").Value = 0.2857142857142857
$ws.Range("S2").Value = 0.6666666666666666
$ws.Range("W2").Value = 0.3333333333333333
$ws.Range("C3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("I3").Value = 0.5
$ws.Range("M3").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("S3").Value = 1
$ws.Range("U3").Value = 0.5
$ws.Range("W3").Value = 0.5
$ws.Range("C4").Value = 0.6666666666666666
$ws.Range("F4").Value = 0.5
$ws.Range("I4").Value = 0.4
$ws.Range("M4").Value = 0.5714285714285715
$ws.Range("Q4").Value = 0.4444444444444445
$ws.Range("S4").Value = 0.8
$ws.Range("U4").Value = 0.6666666666666666
$ws.Range("W4").Value = 0.4
$ws.Range("C5").Value = 0.8333333333333334
$ws.Range("F5").Value = 0.7142857142857143
$ws.Range("I5").Value = 0.4545454545454545
$ws.Range("M5").Value = 0.7692307692307692
$ws.Range("Q5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.9090909090909091
$ws.Range("U5").Value = 0.5555555555555556
$ws.Range("W5").Value = 0.4545454545454545
$ws.Range("C6").Value = 0.882808018370203
$ws.Range("F6").Value = 0.7895959410076381
$ws.Range("I6").Value = 0.319393943239799
$ws.Range("M6").Value = 0.9467676761267002
$ws.Range("Q6").Value = 0.9467676761267002
$ws.Range("S6").Value = 0.6920202103528978
$ws.Range("W6").Value = 0.319393943239799
$ws.Range("C11").Value = 2
$ws.Range("F11").Value = 3
$ws.Range("I11").Value = 1.5
$ws.Range("M11").Value = 2.5
$ws.Range("Q11").Value = 3.5
$ws.Range("S11").Value = 1.5
$ws.Range("U11").Value = 0.5
$ws.Range("W11").Value = 1.5
